# Refresh the cryptos price list (Price + Volume(1h) columns), and reorder a
# few rows whose rank flipped, per the Mar 11 2024 GitHub Actions data pull.
#
# Price cells (column D) are stored as TEXT in the sheet (e.g. "521.51",
# "0.0000320") even though they look numeric. Excel's COM layer auto-detects
# numeric-looking strings and stores them as numbers, which would silently
# drop things like trailing zeros ("86.80" -> 86.8) or render tiny values in
# scientific notation. Prefixing those assignments with a leading apostrophe
# forces a literal-text entry, matching the original cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($ref, $value) {
    # Leading apostrophe = Excel "treat as text" quote-prefix; stops the
    # engine from re-typing numeric-looking strings as numbers/dates.
    $ws.Range($ref).Value = "'" + $value
}

# Row 2 - Bitcoin
Set-Text 'D2' '68.526.79'
$ws.Range('E2').Value = '  -1.27%  '

# Row 3 - Ethereum
Set-Text 'D3' '3.855.18'
$ws.Range('E3').Value = '  -2.28%  '

# Row 4 - TetherUSD (price unchanged)
$ws.Range('E4').Value = '  +0.10%  '

# Row 5 - BNB
Set-Text 'D5' '521.41'
$ws.Range('E5').Value = '  +5.97%  '

# Row 6 - Solana
Set-Text 'D6' '140.64'
$ws.Range('E6').Value = '  -4.42%  '

# Row 7 - XRP (price unchanged)
$ws.Range('E7').Value = '  -2.74%  '

# Row 8 - USDC
Set-Text 'D8' '0.999'
$ws.Range('E8').Value = '  +0.15%  '

# Row 9 - Cardano
Set-Text 'D9' '0.711'
$ws.Range('E9').Value = '  -3.57%  '

# Row 10 - Dogecoin
Set-Text 'D10' '0.167'
$ws.Range('E10').Value = '  -6.21%  '

# Row 11 - ShibaInu
Set-Text 'D11' '0.0000321'
$ws.Range('E11').Value = '  -7.87%  '

# Row 12 - Avalanche
Set-Text 'D12' '41.51'
$ws.Range('E12').Value = '  -3.82%  '

# Row 13 - Polkadot
Set-Text 'D13' '10.43'
$ws.Range('E13').Value = '  -0.37%  '

# Row 14 - WrappedliquidstakedEther2.0
Set-Text 'D14' '4.472.19'
$ws.Range('E14').Value = '  -2.22%  '

# Rows 15/16 swapped rank: Chainlink <-> WrappedEther
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-Text 'D15' '3.899.69'
$ws.Range('E15').Value = '  -1.57%  '

$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-Text 'D16' '21.31'
$ws.Range('E16').Value = '  +6.87%  '

# Row 17 - Uniswap
Set-Text 'D17' '14.07'
$ws.Range('E17').Value = '  -1.55%  '

# Row 18 - TRON (price unchanged)
$ws.Range('E18').Value = '  -2.17%  '

# Row 19 - Polygon (price unchanged)
$ws.Range('E19').Value = '  +2.12%  '

# Row 20 - WrappedBTC
Set-Text 'D20' '68.571.29'
$ws.Range('E20').Value = '  -1.29%  '

# Row 21 - BitcoinCash
Set-Text 'D21' '415.47'
$ws.Range('E21').Value = '  -5.75%  '

# Row 22 - ImmutableX (price unchanged)
$ws.Range('E22').Value = '  +0.64%  '

# Row 23 - InternetComputer(DFINITY)
Set-Text 'D23' '13.97'
$ws.Range('E23').Value = '  -3.86%  '

# Row 24 - Litecoin
Set-Text 'D24' '86.80'
$ws.Range('E24').Value = '  -2.89%  '

# Row 25 - PancakeSwap (price unchanged)
$ws.Range('E25').Value = '  +6.24%  '

# Row 26 - RenderToken
Set-Text 'D26' '11.66'
$ws.Range('E26').Value = '  -2.80%  '

# Row 27 - Filecoin
Set-Text 'D27' '10.47'
$ws.Range('E27').Value = '  -5.94%  '

# Row 28 - EthereumClassic
Set-Text 'D28' '35.46'
$ws.Range('E28').Value = '  -4.88%  '

# Row 29 - Cosmos
Set-Text 'D29' '13.42'
$ws.Range('E29').Value = '  -0.72%  '

# Row 30 - Bittensor
Set-Text 'D30' '678.52'
$ws.Range('E30').Value = '  -4.11%  '

# Row 31 - Hedera
Set-Text 'D31' '0.124'
$ws.Range('E31').Value = '  -5.43%  '

# Row 32 - NEARProtocol
Set-Text 'D32' '6.80'
$ws.Range('E32').Value = '  +12.09%  '

# Row 33 - Toncoin (price unchanged)
$ws.Range('E33').Value = '  -3.71%  '

# Row 34 - OKB
Set-Text 'D34' '66.97'
$ws.Range('E34').Value = '  +8.66%  '

# Row 35 - TheGraph
Set-Text 'D35' '0.443'
$ws.Range('E35').Value = '  -6.85%  '

# Rows 36/37 swapped rank: InjectiveProtocol <-> PEPE
$ws.Range('B36').Value = 'PEPE'
$ws.Range('C36').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D36').Value = '0.0₃0839'
$ws.Range('E36').Value = '  -8.23%  '

$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-Text 'D37' '39.40'
$ws.Range('E37').Value = '  -3.60%  '

# Row 38 - ThetaToken
Set-Text 'D38' '3.48'
$ws.Range('E38').Value = '  +12.97%  '

# Row 39 - Kaspa
Set-Text 'D39' '0.148'
$ws.Range('E39').Value = '  -1.18%  '

# Row 40 - Dai (price unchanged)
$ws.Range('E40').Value = '  +0.07%  '

# Row 41 - FirstDigitalUSD (price unchanged)
$ws.Range('E41').Value = '  -0.23%  '

# Row 42 - WEMIXToken (price unchanged)
$ws.Range('E42').Value = '  +5.35%  '

# Row 43 - VeChain
Set-Text 'D43' '0.0473'
$ws.Range('E43').Value = '  -3.64%  '

# Row 44 - Fetch.AI
Set-Text 'D44' '2.82'
$ws.Range('E44').Value = '  -3.99%  '

# Row 45 - ApeXProtocol
Set-Text 'D45' '3.41'
$ws.Range('E45').Value = '  +1.95%  '

# Row 46 - Stellar (price unchanged)
$ws.Range('E46').Value = '  -1.84%  '

# Row 47 - Stacks
Set-Text 'D47' '2.98'
$ws.Range('E47').Value = '  -2.86%  '

# Row 48 - FLOKI
Set-Text 'D48' '0.000271'
$ws.Range('E48').Value = '  +11.32%  '

# Rows 49/50 swapped rank: LidoDAOToken <-> Monero
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-Text 'D49' '143.84'
$ws.Range('E49').Value = '  -0.10%  '

$ws.Range('B50').Value = 'LidoDAOToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-Text 'D50' '3.26'
$ws.Range('E50').Value = '  -3.50%  '

# Row 51 - BabyDogeCoin (price unchanged)
$ws.Range('E51').Value = '  -8.03%  '
